# Insert a new weekly price record as row 121 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 121-186 down to 122-187 (the former last row,
# 186, simply moves to 187 unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(121).Insert()

$ws.Cells.Item(121, 1).Value  = 4
$ws.Cells.Item(121, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(121, 3).Value  = "Los Lagos"
$ws.Cells.Item(121, 4).Value  = 44460
$ws.Cells.Item(121, 5).Value  = 10
$ws.Cells.Item(121, 6).Value  = 100112023
$ws.Cells.Item(121, 7).Value  = "Brócoli"
$ws.Cells.Item(121, 8).Value  = "Sin especificar"
$ws.Cells.Item(121, 9).Value  = "Segunda"
$ws.Cells.Item(121, 10).Value = 500
$ws.Cells.Item(121, 11).Value = 1000
$ws.Cells.Item(121, 12).Value = 1000
$ws.Cells.Item(121, 13).Value = 1000
$ws.Cells.Item(121, 14).Value = "$/unidad"
$ws.Cells.Item(121, 15).Value = "Región del Maule"
$ws.Cells.Item(121, 16).Value = 1000
$ws.Cells.Item(121, 17).Value = 1
$ws.Cells.Item(121, 18).Value = "Hortaliza"
